$wb = $excel.ActiveWorkbook

# Update scraped_at timestamps on the "snapshot" sheet (K2:K50)
$snapshot = $wb.Worksheets.Item("snapshot")

$snapshot.Range("K2").Value = "2025-10-30T03:01:40.962534+00:00"
$snapshot.Range("K3").Value = "2025-10-30T03:01:43.761920+00:00"
$snapshot.Range("K4").Value = "2025-10-30T03:01:43.761950+00:00"
$snapshot.Range("K5").Value = "2025-10-30T03:01:43.761969+00:00"
$snapshot.Range("K6").Value = "2025-10-30T03:01:43.761986+00:00"
$snapshot.Range("K7").Value = "2025-10-30T03:01:43.762002+00:00"
$snapshot.Range("K8").Value = "2025-10-30T03:01:46.084856+00:00"
$snapshot.Range("K9").Value = "2025-10-30T03:01:46.084884+00:00"
$snapshot.Range("K10").Value = "2025-10-30T03:01:48.843727+00:00"
$snapshot.Range("K11").Value = "2025-10-30T03:01:51.634871+00:00"
$snapshot.Range("K12").Value = "2025-10-30T03:01:51.634904+00:00"
$snapshot.Range("K13").Value = "2025-10-30T03:01:53.926844+00:00"
$snapshot.Range("K14").Value = "2025-10-30T03:01:53.926875+00:00"
$snapshot.Range("K15").Value = "2025-10-30T03:01:53.926891+00:00"
$snapshot.Range("K16").Value = "2025-10-30T03:01:53.926907+00:00"
$snapshot.Range("K17").Value = "2025-10-30T03:02:01.663220+00:00"
$snapshot.Range("K18").Value = "2025-10-30T03:02:04.428959+00:00"
$snapshot.Range("K19").Value = "2025-10-30T03:02:07.293863+00:00"
$snapshot.Range("K20").Value = "2025-10-30T03:02:10.076586+00:00"
$snapshot.Range("K21").Value = "2025-10-30T03:02:10.076625+00:00"
$snapshot.Range("K22").Value = "2025-10-30T03:02:10.076647+00:00"
$snapshot.Range("K23").Value = "2025-10-30T03:02:12.760114+00:00"
$snapshot.Range("K24").Value = "2025-10-30T03:02:12.760141+00:00"
$snapshot.Range("K25").Value = "2025-10-30T03:02:12.760158+00:00"
$snapshot.Range("K26").Value = "2025-10-30T03:02:12.760173+00:00"
$snapshot.Range("K27").Value = "2025-10-30T03:02:12.760189+00:00"
$snapshot.Range("K28").Value = "2025-10-30T03:02:20.817076+00:00"
$snapshot.Range("K29").Value = "2025-10-30T03:02:20.817114+00:00"
$snapshot.Range("K30").Value = "2025-10-30T03:02:20.817138+00:00"
$snapshot.Range("K31").Value = "2025-10-30T03:02:20.817159+00:00"
$snapshot.Range("K32").Value = "2025-10-30T03:02:23.436599+00:00"
$snapshot.Range("K33").Value = "2025-10-30T03:02:23.436627+00:00"
$snapshot.Range("K34").Value = "2025-10-30T03:02:23.436645+00:00"
$snapshot.Range("K35").Value = "2025-10-30T03:02:25.644733+00:00"
$snapshot.Range("K36").Value = "2025-10-30T03:02:25.644761+00:00"
$snapshot.Range("K37").Value = "2025-10-30T03:02:25.644778+00:00"
$snapshot.Range("K38").Value = "2025-10-30T03:02:25.644793+00:00"
$snapshot.Range("K39").Value = "2025-10-30T03:02:25.644808+00:00"
$snapshot.Range("K40").Value = "2025-10-30T03:02:25.644822+00:00"
$snapshot.Range("K41").Value = "2025-10-30T03:02:25.644837+00:00"
$snapshot.Range("K42").Value = "2025-10-30T03:02:25.644857+00:00"
$snapshot.Range("K43").Value = "2025-10-30T03:02:25.644871+00:00"
$snapshot.Range("K44").Value = "2025-10-30T03:02:28.430295+00:00"
$snapshot.Range("K45").Value = "2025-10-30T03:02:28.430372+00:00"
$snapshot.Range("K46").Value = "2025-10-30T03:02:33.979146+00:00"
$snapshot.Range("K47").Value = "2025-10-30T03:02:36.779689+00:00"
$snapshot.Range("K48").Value = "2025-10-30T03:02:36.779717+00:00"
$snapshot.Range("K49").Value = "2025-10-30T03:02:36.779734+00:00"
$snapshot.Range("K50").Value = "2025-10-30T03:02:36.779750+00:00"

# Remove the returned player row (row 2) from the "returned" sheet,
# leaving only the header row (dimension becomes A1:G1)
$returned = $wb.Worksheets.Item("returned")
$returned.Rows.Item(2).Delete()

